# Updates the "sps" worksheet's Jinja template layout for single search parameters.
# - Collapses the single modifier/comparator/chain columns into SHALL/SHOULD pairs.
# - Moves "references" from column Q to column S.
# - Refreshes example chain/comparator example data in rows 5, 6, 11, 16.
# - Updates the sheet view's scroll position / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sps")

# ---- Row 1 headers (columns M:S) ----
$ws.Range("M1").Value = "shall_modifier"
$ws.Range("N1").Value = "should_modifier"
$ws.Range("O1").Value = "shall_comparator"
$ws.Range("P1").Value = "should_comparator"
$ws.Range("Q1").Value = "shall_chain"
$ws.Range("R1").Value = "should_chain"
$ws.Range("S1").Value = "references"

# ---- Clear the old per-row "n/a" filler data in columns M:S (rows 2-16) ----
$ws.Range("M2:S16").ClearContents()

# ---- Re-populate the handful of rows that still carry data in the new layout ----

# Row 5 (Questionnaire.title) - SHALL/SHOULD comparator + chain examples
$ws.Range("N5").Value = "foo,bar2"
$ws.Range("O5").Value = "contains,foo,bar"
$ws.Range("P5").Value = "foo,bar2"
$ws.Range("Q5").Value = "contains,foo,bar"
$ws.Range("R5").Value = "foo,bar2"

# Row 6 (Questionnaire.publisher) - shall_modifier stays "contains"
$ws.Range("M6").Value = "contains"

# Row 11 (QuestionnaireResponse.patient) - shall_chain / should_chain targets
$ws.Range("Q11").Value = "patient.gender,patient.birthdate,patient.race,patient.ethnicity,patient.address-postalcode"
$ws.Range("R11").Value = "SHALL,SHALL,SHALL,SHALL,SHALL"

# Row 16 (Patient.birthdate) - shall_comparator stays "le,ge"
$ws.Range("O16").Value = "le,ge"

# ---- Sheet view: scroll position + active selection ----
$ws.Range("P5").Select()
